$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "69.310.14"
$ws.Range("E2").Value = "  -0.04%  "

# Row 3 - Ethereum
Set-TextCell "D3" "3.675.49"
$ws.Range("E3").Value = "  -0.30%  "

# Row 4 - TetherUSD
Set-TextCell "D4" "1.00"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
Set-TextCell "D5" "682.60"
$ws.Range("E5").Value = "  -0.31%  "

# Row 6 - Solana
Set-TextCell "D6" "157.93"
$ws.Range("E6").Value = "  -2.93%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -1.42%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -2.02%  "

# Row 10 - Toncoin
Set-TextCell "D10" "6.98"
$ws.Range("E10").Value = "  -3.81%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -3.52%  "

# Row 12 - ShibaInu
$ws.Range("E12").Value = "  -2.21%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextCell "D13" "4.295.30"
$ws.Range("E13").Value = "  -0.33%  "

# Row 14 - Avalanche
Set-TextCell "D14" "32.12"
$ws.Range("E14").Value = "  -4.50%  "

# Row 15 - WrappedEther
Set-TextCell "D15" "3.663.19"
$ws.Range("E15").Value = "  -0.59%  "

# Row 16 - WrappedBTC
Set-TextCell "D16" "69.299.12"
$ws.Range("E16").Value = "  -0.13%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  +2.02%  "

# Row 18 - Chainlink
Set-TextCell "D18" "15.85"
$ws.Range("E18").Value = "  -2.75%  "

# Row 19 - Polkadot
Set-TextCell "D19" "6.38"
$ws.Range("E19").Value = "  -4.04%  "

# Row 20 - BitcoinCash
Set-TextCell "D20" "471.55"
$ws.Range("E20").Value = "  -2.03%  "

# Row 21 - Uniswap
Set-TextCell "D21" "9.94"
$ws.Range("E21").Value = "  +1.19%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  -3.04%  "

# Row 23 - Litecoin
$ws.Range("E23").Value = "  -0.14%  "

# Row 24 - WrappedeETH
Set-TextCell "D24" "3.820.49"
$ws.Range("E24").Value = "  -0.30%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  -0.10%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  -5.70%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("E27").Value = "  -5.15%  "

# Row 28 - RenderToken
$ws.Range("E28").Value = "  -4.91%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  -2.05%  "

# Row 30 - Fetch.AI
Set-TextCell "D30" "1.73"
$ws.Range("E30").Value = "  -5.65%  "

# Row 31 - Binance-PegBSC-USD
$ws.Range("E31").Value = "  +0.09%  "

# Row 32 - NEARProtocol
Set-TextCell "D32" "6.53"
$ws.Range("E32").Value = "  -4.46%  "

# Row 33 - ImmutableX
$ws.Range("E33").Value = "  -6.37%  "

# Row 34 - EthereumClassic
$ws.Range("E34").Value = "  -1.05%  "

# Row 35 - RenzoRestakedETH
Set-TextCell "D35" "3.653.46"
$ws.Range("E35").Value = "  +0.11%  "

# Row 36 - Kaspa
$ws.Range("E36").Value = "  -4.09%  "

# Row 37 - Aptos
Set-TextCell "D37" "8.15"
$ws.Range("E37").Value = "  -5.37%  "

# Row 38 - Filecoin
Set-TextCell "D38" "6.07"
$ws.Range("E38").Value = "  -0.27%  "

# Row 40 - Stacks
Set-TextCell "D40" "2.21"
$ws.Range("E40").Value = "  +1.51%  "

# Row 41 - Hedera
Set-TextCell "D41" "0.0897"

# Row 42 - FirstDigitalUSD
$ws.Range("E42").Value = "  -0.07%  "

# Row 43 - Mantle
$ws.Range("E43").Value = "  -2.11%  "

# Row 44 - Monero
Set-TextCell "D44" "166.20"
$ws.Range("E44").Value = "  +5.24%  "

# Row 45 - OKB
Set-TextCell "D45" "47.56"
$ws.Range("E45").Value = "  -1.12%  "

# Row 46 - now dogwifhat (was FLOKI)
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell "D46" "2.71"
$ws.Range("E46").Value = "  -5.06%  "

# Row 47 - now FLOKI (was dogwifhat)
$ws.Range("B47").Value = "FLOKI"
$ws.Range("C47").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextCell "D47" "0.000277"
$ws.Range("E47").Value = "  -1.71%  "

# Row 48 - SuiNetwork
Set-TextCell "D48" "1.10"
$ws.Range("E48").Value = "  +0.86%  "

# Row 49 - ONDO
$ws.Range("E49").Value = "  -3.00%  "

# Row 50 - Cosmos
$ws.Range("E50").Value = "  -4.52%  "

# Row 51 - InjectiveProtocol
Set-TextCell "D51" "26.74"
$ws.Range("E51").Value = "  -4.41%  "
